{"js": "// Edit 1: \"Para la refrigeraci\u00f3n del edificio hemos...\" ->\n//          \"Para la refrigeraci\u00f3n del cuerpo principal del edificio hemos...\"\n// We locate the unique anchor \"edificio hemos\" and insert the extra wording\n// right before it so the surrounding text matches exactly.\nconst anchor1 = context.document.body.search(\"edificio hemos\", { matchCase: true, matchWholeWord: false });\nanchor1.load(\"items\");\nawait context.sync();\n\nif (anchor1.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the first anchor, found \" + anchor1.items.length);\n}\nanchor1.items[0].insertText(\"cuerpo principal del \", Word.InsertLocation.before);\nawait context.sync();\n\n// Edit 2: add a new paragraph right after the \"Para la habitaci\u00f3n de\n// servidores...\" paragraph, and move the _GoBack bookmark range (the end\n// of the document body) onto the new paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet serverParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Para la habitaci\u00f3n de servidores\") !== -1) {\n    serverParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!serverParagraph) {\n  throw new Error(\"Could not find the servers paragraph\");\n}\n\nserverParagraph.insertParagraph(\n  \"Adem\u00e1s habr\u00e1 dos aire acondicionado de 18 mil frigor\u00edas en el  edificio de expedici\u00f3n/producci\u00f3n, que funcionar\u00e1n en la noche, horario en que trabajan las rotativas.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Re-load the paragraph collection so we get a \"live\" reference to the\n// freshly-inserted paragraph (the object returned directly by\n// insertParagraph does not reliably resolve further range queries).\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\nconst newParagraph = refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\n\n// The \"_GoBack\" bookmark previously sat at the end of the servers\n// paragraph; move it to the end of the newly-added paragraph (matching\n// where Word leaves it after the last edit location).\ncontext.document.deleteBookmark(\"_GoBack\");\nconst newParagraphEnd = newParagraph.getRange(\"End\");\nnewParagraphEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1 -----------------------------------------------------------\n# \"Para la refrigeraci\u00f3n del edificio hemos...\" ->\n# \"Para la refrigeraci\u00f3n del cuerpo principal del edificio hemos...\"\n# \"edificio hemos\" is unique in the document, so find it and insert the\n# extra wording right before it.\n$findRange = $d.Content\n$findRange.Find.MatchCase = $true\n$found = $findRange.Find.Execute(\"edificio hemos\")\nif (-not $found) {\n    throw \"Could not find the anchor text 'edificio hemos'\"\n}\n$findRange.InsertBefore(\"cuerpo principal del \")\n\n# --- Edit 2 -------------------------------------------------------------\n# Add a new paragraph right after the \"Para la habitaci\u00f3n de\n# servidores...\" paragraph, and move the hidden \"_GoBack\" bookmark so it\n# again sits at the very end of the document (now the end of the new\n# paragraph instead of the end of the servers paragraph).\n\n# Locate the servers paragraph index by its distinctive text (re-reading\n# the paragraph by index below rather than holding on to the COM object\n# keeps every access fresh/live).\n$serverIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text -like \"*Para la habitaci*n de servidores*\") {\n        $serverIndex = $i\n        break\n    }\n}\nif ($serverIndex -eq -1) {\n    throw \"Could not find the servers paragraph\"\n}\n\n$newParagraphText = \"Adem\u00e1s habr\u00e1 dos aire acondicionado de 18 mil frigor\u00edas en el  edificio de expedici\u00f3n/producci\u00f3n, que funcionar\u00e1n en la noche, horario en que trabajan las rotativas.\"\n\n$d.Paragraphs.Item($serverIndex).Range.InsertParagraphAfter()\n$newIndex = $serverIndex + 1\n\n# Set the new paragraph's text with a one-character sentinel suffix. While\n# we create the \"_GoBack\" bookmark, its insertion point must never be the\n# very last character position in the document (that specific edge case\n# does not reliably place the bookmark). We add the bookmark just before\n# the sentinel, then delete the sentinel, leaving a correctly-collapsed\n# bookmark at the true end of the paragraph.\n$d.Paragraphs.Item($newIndex).Range.Text = $newParagraphText + \"X\"\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$sentinelPos = $d.Paragraphs.Item($newIndex).Range.Start + $newParagraphText.Length\n$bookmarkPoint = $d.Range($sentinelPos, $sentinelPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n\n$sentinelEnd = $d.Paragraphs.Item($newIndex).Range.End - 1\n$sentinelRange = $d.Range($sentinelPos, $sentinelEnd)\n$sentinelRange.Text = \"\"\n"}
